$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.277.36"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "1.819.63"
$ws.Range("E3").Value = "  +1.25%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").Value = "'313.08"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.4660"
$ws.Range("E7").Value = "  +4.48%  "

$ws.Range("D8").Value = "'0.3774"
$ws.Range("E8").Value = "  +2.47%  "

$ws.Range("D9").Value = "'0.07405"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("D10").Value = "'0.8700"
$ws.Range("E10").Value = "  +1.25%  "

$ws.Range("D11").Value = "'20.62"
$ws.Range("E11").Value = "  -0.14%  "

$ws.Range("D12").Value = "1.823.63"
$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("D13").Value = "'6.685"
$ws.Range("E13").Value = "  +0.88%  "

$ws.Range("D14").Value = "'5.411"
$ws.Range("E14").Value = "  +2.83%  "

$ws.Range("D15").Value = "'92.30"
$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("D16").Value = "'0.07087"
$ws.Range("E16").Value = "  +0.27%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").Value = "'0.000008760"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("D21").Value = "27.271.59"
$ws.Range("E21").Value = "  +1.95%  "

$ws.Range("D22").Value = "'5.313"
$ws.Range("E22").Value = "  +3.01%  "

$ws.Range("D23").Value = "'10.92"
$ws.Range("E23").Value = "  +1.31%  "

$ws.Range("D24").Value = "2.049.27"
$ws.Range("E24").Value = "  +1.45%  "

$ws.Range("D25").Value = "'1.940"
$ws.Range("E25").Value = "  -2.38%  "

$ws.Range("D26").Value = "'151.48"
$ws.Range("E26").Value = "  -0.22%  "

$ws.Range("D27").Value = "'2.245"
$ws.Range("E27").Value = "  +3.62%  "

$ws.Range("D28").Value = "'18.53"
$ws.Range("E28").Value = "  +0.79%  "

$ws.Range("D29").Value = "'5.315"
$ws.Range("E29").Value = "  +2.35%  "

$ws.Range("D30").Value = "'117.06"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("D31").Value = "'0.08932"
$ws.Range("E31").Value = "  +1.72%  "

$ws.Range("D32").Value = "'0.7820"
$ws.Range("E32").Value = "  +5.75%  "

$ws.Range("E33").Value = "  +2.36%  "

$ws.Range("D34").Value = "'4.526"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").Value = "'2.926"
$ws.Range("E35").Value = "  +0.29%  "

$ws.Range("D36").Value = "'0.9999"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").Value = "'1.097"
$ws.Range("E37").Value = "  +1.37%  "

$ws.Range("D38").Value = "'0.01966"
$ws.Range("E38").Value = "  +0.50%  "

$ws.Range("D39").Value = "'0.05249"
$ws.Range("E39").Value = "  +1.19%  "

$ws.Range("D40").Value = "'7.279"
$ws.Range("E40").Value = "  +4.60%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5312"
$ws.Range("E41").Value = "  +0.83%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'2.370"
$ws.Range("E42").Value = "  +20.00%  "

$ws.Range("D43").Value = "'2.890"
$ws.Range("E43").Value = "  +2.07%  "

$ws.Range("D44").Value = "'0.1690"
$ws.Range("E44").Value = "  +0.56%  "

$ws.Range("D45").Value = "'8.612"
$ws.Range("E45").Value = "  +2.28%  "

$ws.Range("D46").Value = "'0.5056"
$ws.Range("E46").Value = "  -0.49%  "

$ws.Range("D47").Value = "'10.44"
$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").Value = "'105.58"
$ws.Range("E48").Value = "  +0.58%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").Value = "'1.669"
$ws.Range("E50").Value = "  -0.08%  "

$ws.Range("D51").Value = "'0.06328"
$ws.Range("E51").Value = "  +0.57%  "
